# Swap the deck's theme palette: the single Design ("Integral") that backs
# every slide (ppt/theme/theme1.xml, reached through SlideMaster/slide
# ThemeColorScheme) switches from the Integral color scheme to the stock
# Office Theme color scheme.
#
# PowerPoint's VBA RGB() macro packs r,g,b into a little-endian OLE COLORREF
# (0x00BBGGRR) -- reproduce that here so ThemeColor.RGB gets the right value.
function RGB($r, $g, $b) {
    return [int]$r -bor ([int]$g * 256) -bor ([int]$b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme color scheme (the standard default PowerPoint theme palette)
# in the fixed ThemeColorScheme index order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

# The ThemeColorScheme is shared by the whole deck (it is backed by the
# slide master's theme part), so reaching it through the first slide
# repaints every slide + the slide master in one go.
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
